# Auto-committed on 2022/02/25 週五
#
# Change the column-D "type" values for the AuthItem (row 10) and Desc
# (row 11) fields on the DBD sheet from "varchar2" to "nvarchar2", and
# switch the active sheet/selection from DBS (B5) back to DBD (G10).

$wb  = $excel.ActiveWorkbook
$dbd = $wb.Sheets.Item("DBD")
$dbs = $wb.Sheets.Item("DBS")

# AuthItem / Desc 型態 varchar2 -> nvarchar2
$dbd.Range("D10").Value = "nvarchar2"
$dbd.Range("D11").Value = "nvarchar2"

# Make DBD the active sheet/tab again, with G10 selected.
$dbd.Activate()
$dbd.Range("G10").Select()
